$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 6 new columns at their final positions, left to right.
# Order matters for shared-string allocation order (matches original author's edit order).
$ws.Columns("A").Insert()
$ws.Columns("G").Insert()
$ws.Columns("H").Insert()
$ws.Columns("K").Insert()
$ws.Columns("M").Insert()
$ws.Columns("O").Insert()

# Fill in the header values for the newly inserted columns, in the same
# order the columns were created so shared strings end up in the expected
# sequence.
$ws.Range("G1").Value = "RENEWED_UPTO"
$ws.Range("H1").Value = "SCHOOL_STATUS"
$ws.Range("A1").Value = "CLASS"
$ws.Range("K1").Value = "CANDIDATE_ADHAAR"
$ws.Range("M1").Value = "MOTHER_ADHAAR"
$ws.Range("O1").Value = "FATHER_ADHAAR"

# Column widths for the new columns (not bestFit, explicit custom widths).
$ws.Range("G1:H1").ColumnWidth = 15.166666666666666
$ws.Range("K1").ColumnWidth = 17.166666666666668
$ws.Range("M1").ColumnWidth = 16.5
$ws.Range("O1").ColumnWidth = 16

# Text number format on the two "DOB"/"date cutting" style columns (now Q and AA).
$ws.Range("Q1").NumberFormat = "@"
$ws.Range("AA1").NumberFormat = "@"

# Page setup.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update selection / scroll position.
$ws.Range("E11").Select()
